$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 323, shifting existing rows 323+ down by one.
$ws.Rows.Item(323).Insert()

# Populate the new row 323 with the new data values.
$ws.Cells.Item(323, 1).Value = 6
$ws.Cells.Item(323, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(323, 3).Value = "Metropolitana"
$ws.Cells.Item(323, 4).Value = 44505
$ws.Cells.Item(323, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(323, 5).Value = 13
$ws.Cells.Item(323, 6).Value = 100112044
$ws.Cells.Item(323, 7).Value = "Perejil"
$ws.Cells.Item(323, 8).Value = "Sin especificar"
$ws.Cells.Item(323, 9).Value = "Primera"
$ws.Cells.Item(323, 10).Value = 220
$ws.Cells.Item(323, 11).Value = 12000
$ws.Cells.Item(323, 12).Value = 13000
$ws.Cells.Item(323, 13).Value = 12591
$ws.Cells.Item(323, 14).Value = "$/docena de atados"
$ws.Cells.Item(323, 15).Value = "Región Metropolitana"
$ws.Cells.Item(323, 16).Value = 4197
$ws.Cells.Item(323, 17).Value = 3
$ws.Cells.Item(323, 18).Value = "Hortaliza"
